# Add "Service Center Poc Phone" (Mobile No) column to the Spare Requested
# Parts template, between "Service Center Name" (I) and "Service Center
# Address" (old J, now K). CRM-1338.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J; it inherits formatting from the column before it
# (same as Excel's native "Insert Column" behaviour), so the header/placeholder
# rows keep their look automatically.
$ws.Columns.Item(10).Insert()

$ws.Range("J1").Value = "Service Center Mobile No"
$ws.Range("J2").Value = "{spare:primary_contact_phone_1}"

# Match the width of the neighbouring "Service Center Name" column.
$ws.Columns.Item(10).ColumnWidth = $ws.Columns.Item(9).ColumnWidth

# The (now shifted) "Service Center GST Number" header picks up a bold weight
# as part of this change.
$ws.Range("O1").Font.Bold = $true
